# bổ sung chú thích
# Replace the "Kết quả" (result) column values: "có"/"không" -> "nghỉ"/"học"
# per the updated annotation scheme for the lớp học (class) dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "nghỉ"
$ws.Range("C3").Value = "nghỉ"
$ws.Range("C4").Value = "nghỉ"
$ws.Range("C5").Value = "học"
$ws.Range("C6").Value = "học"
$ws.Range("C7").Value = "học"
$ws.Range("C8").Value = "nghỉ"

# Match the author's final selection (cell C8 active).
$ws.Range("C8").Select()
